# Update the embedded build timestamp throughout the workbook.
# Old timestamp: February 03 2026 17.29.55 EST
# New timestamp: February 03 2026 18.05.36 EST

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    foreach ($cell in $usedRange.Cells) {
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
            $cell.Value = $val.Replace($oldStamp, $newStamp)
        }
    }
}
